$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 115.2
$ws.Range("J9").Value = 42.5
$ws.Range("L9").Value = 42.5
$ws.Range("N9").Value = -380.5
$ws.Range("H33").Value = 89
$ws.Range("I33").Value = 89
$ws.Range("K33").Value = 89
$ws.Range("M33").Value = 140
$ws.Range("H40").Value = 2499.7273
$ws.Range("I40").Value = 2199.4
$ws.Range("J40").Value = 2750
$ws.Range("K40").Value = 2199.4
$ws.Range("L40").Value = 2750
$ws.Range("M40").Value = -2024.4
$ws.Range("N40").Value = -3100
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H86").Value = 3146.125
$ws.Range("J86").Value = 4000
$ws.Range("L86").Value = 4000
$ws.Range("N86").Value = -6246
$ws.Range("H89").Value = 3146.125
$ws.Range("J89").Value = 4000
$ws.Range("L89").Value = 20000
$ws.Range("N89").Value = -31232
$ws.Range("H98").Value = 2011.8572
$ws.Range("I98").Value = 2011.8572
$ws.Range("K98").Value = 2011.8572
$ws.Range("M98").Value = -513.8571999999999
$ws.Range("H112").Value = 1956.7
$ws.Range("J112").Value = 2227.125
$ws.Range("L112").Value = 6681.375
$ws.Range("N112").Value = -8897.375
$ws.Range("H122").Value = 2011.8572
$ws.Range("I122").Value = 2011.8572
$ws.Range("K122").Value = 6035.571599999999
$ws.Range("M122").Value = -3585.571599999999
$ws.Range("H138").Value = 3950.2856
$ws.Range("I138").Value = 3171.9697
$ws.Range("K138").Value = 9515.909100000001
$ws.Range("M138").Value = -4375.909100000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 1933.1666
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("H32").Value = 6900.2144
$ws.Range("I32").Value = 4722.6113
$ws.Range("K32").Value = 4722.6113
$ws.Range("M32").Value = -4435.6113
$ws.Range("H92").Value = 70000
$ws.Range("J92").Value = 70000
$ws.Range("L92").Value = 70000
$ws.Range("N92").Value = -74992
$ws.Range("H110").Value = 7239.4443
$ws.Range("I110").Value = 8519.4
$ws.Range("K110").Value = 8519.4
$ws.Range("M110").Value = -6474.4
$ws.Range("H138").Value = 74530
$ws.Range("J138").Value = 74530
$ws.Range("L138").Value = 74530
$ws.Range("N138").Value = -84810
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 13824
$ws.Range("I26").Value = 13824
$ws.Range("K26").Value = 13824
$ws.Range("M26").Value = -13532
$ws.Range("H88").Value = 31233.715
$ws.Range("I88").Value = 15300
$ws.Range("J88").Value = 33889.332
$ws.Range("K88").Value = 15300
$ws.Range("L88").Value = 33889.332
$ws.Range("M88").Value = -14894
$ws.Range("N88").Value = -34701.332
$ws.Range("H91").Value = 31233.715
$ws.Range("I91").Value = 15300
$ws.Range("J91").Value = 33889.332
$ws.Range("K91").Value = 15300
$ws.Range("L91").Value = 33889.332
$ws.Range("M91").Value = -13896
$ws.Range("N91").Value = -36697.332
$ws.Range("H99").Value = 2852.8667
$ws.Range("J99").Value = 2800
$ws.Range("L99").Value = 2800
$ws.Range("N99").Value = -5796
$ws.Range("H105").Value = 3133.6155
$ws.Range("I105").Value = 2885.2727
$ws.Range("K105").Value = 2885.2727
$ws.Range("M105").Value = -1138.2727
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3747
$ws.Range("I16").Value = 3746.5
$ws.Range("K16").Value = 3746.5
$ws.Range("M16").Value = -3459.5
$ws.Range("H22").Value = 367.52942
$ws.Range("I22").Value = 245.25
$ws.Range("K22").Value = 245.25
$ws.Range("M22").Value = 104.75
$ws.Range("H31").Value = 2874
$ws.Range("I31").Value = 2499.4075
$ws.Range("K31").Value = 2499.4075
$ws.Range("M31").Value = -2204.4075
$ws.Range("H34").Value = 2874
$ws.Range("I34").Value = 2499.4075
$ws.Range("K34").Value = 2499.4075
$ws.Range("M34").Value = -2297.4075
$ws.Range("H86").Value = 15000
$ws.Range("J86").Value = 15000
$ws.Range("L86").Value = 15000
$ws.Range("N86").Value = -17246
$ws.Range("H88").Value = 90258
$ws.Range("J88").Value = 90258
$ws.Range("L88").Value = 90258
$ws.Range("N88").Value = -91070
$ws.Range("H89").Value = 15000
$ws.Range("J89").Value = 15000
$ws.Range("L89").Value = 75000
$ws.Range("N89").Value = -86232
$ws.Range("H91").Value = 90258
$ws.Range("J91").Value = 90258
$ws.Range("L91").Value = 90258
$ws.Range("N91").Value = -93066
$ws.Range("H107").Value = 1170.8889
$ws.Range("I107").Value = 691.6667
$ws.Range("K107").Value = 691.6667
$ws.Range("M107").Value = 1228.3333
$ws.Range("H113").Value = 3747
$ws.Range("I113").Value = 3746.5
$ws.Range("K113").Value = 3746.5
$ws.Range("M113").Value = -1576.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 79.25
$ws.Range("I40").Value = 43.5
$ws.Range("K40").Value = 174
$ws.Range("M40").Value = -105
$ws.Range("H68").Value = 1516.8
$ws.Range("J68").Value = 2142
$ws.Range("L68").Value = 6426
$ws.Range("N68").Value = -8048
$ws.Range("H71").Value = 1516.8
$ws.Range("J71").Value = 2142
$ws.Range("L71").Value = 19278
$ws.Range("N71").Value = -27390
$ws.Range("H92").Value = 317.5
$ws.Range("I92").Value = 330.9091
$ws.Range("J92").Value = 296.42856
$ws.Range("K92").Value = 992.7273
$ws.Range("L92").Value = 889.28568
$ws.Range("M92").Value = 255.2727
$ws.Range("N92").Value = -3385.28568
$ws.Range("H141").Value = 5257.25
$ws.Range("I141").Value = 5257.25
$ws.Range("K141").Value = 15771.75
$ws.Range("M141").Value = -10591.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1033.3334
$ws.Range("I107").Value = 800
$ws.Range("K107").Value = 800
$ws.Range("M107").Value = 1120
$ws.Range("H113").Value = 4250
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 5883438.5
$ws.Range("I13").Value = 6250937.5
$ws.Range("K13").Value = 6250937.5
$ws.Range("M13").Value = -6250797.5
$ws.Range("H46").Value = 4959.25
$ws.Range("I46").Value = 999
$ws.Range("K46").Value = 999
$ws.Range("M46").Value = -811
$ws.Range("H61").Value = 3863
$ws.Range("I61").Value = 3720.875
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 3720.875
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -3518.875
$ws.Range("N61").Value = -5404
$ws.Range("H76").Value = 11500
$ws.Range("J76").Value = 11500
$ws.Range("L76").Value = 11500
$ws.Range("N76").Value = -12176
$ws.Range("H79").Value = 11500
$ws.Range("J79").Value = 11500
$ws.Range("L79").Value = 11500
$ws.Range("N79").Value = -13840
$ws.Range("H113").Value = 3863
$ws.Range("I113").Value = 3720.875
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 3720.875
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -1550.875
$ws.Range("N113").Value = -9340
$ws.Range("H136").Value = 5512
$ws.Range("I136").Value = 3574.5
$ws.Range("J136").Value = 8095.3335
$ws.Range("K136").Value = 10723.5
$ws.Range("L136").Value = 24286.0005
$ws.Range("M136").Value = -8173.5
$ws.Range("N136").Value = -29386.0005
